$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.428.10"
$ws.Range("E2").Value = "  +2.78%  "
$ws.Range("D3").Value = "1.837.68"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.89"
$ws.Range("E5").Value = "  +2.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  +1.50%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.95"
$ws.Range("E8").Value = "  +13.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.310"
$ws.Range("E9").Value = "  +7.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0705"
$ws.Range("E10").Value = "  +5.51%  "
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("D12").Value = "2.104.34"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").Value = "1.838.05"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("E15").Value = "  +6.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.69"
$ws.Range("E16").Value = "  +7.29%  "
$ws.Range("D17").Value = "35.429.18"
$ws.Range("E17").Value = "  +2.79%  "
$ws.Range("E18").Value = "  +2.70%  "
$ws.Range("E19").Value = "  +4.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "244.20"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.10"
$ws.Range("E21").Value = "  +8.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.69"
$ws.Range("E22").Value = "  +14.62%  "
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.12"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.90"
$ws.Range("E26").Value = "  +3.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.74"
$ws.Range("E27").Value = "  +1.38%  "
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("E29").Value = "  +22.07%  "
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").Value = "3.355.41"
$ws.Range("E31").Value = "  +38.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0553"
$ws.Range("E32").Value = "  +7.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.10"
$ws.Range("E33").Value = "  +6.78%  "
$ws.Range("E34").Value = "  +4.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.86"
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "96.06"
$ws.Range("E36").Value = "  +16.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.686"
$ws.Range("E37").Value = "  +7.48%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.83"
$ws.Range("E38").Value = "  +13.10%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.10"
$ws.Range("E39").Value = "  +4.36%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.348.75"
$ws.Range("E40").Value = "  +3.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0195"
$ws.Range("E41").Value = "  +4.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.43"
$ws.Range("E42").Value = "  +5.06%  "
$ws.Range("E43").Value = "  +6.33%  "
$ws.Range("E44").Value = "  +4.48%  "
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.26"
$ws.Range("E47").Value = "  +8.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0520"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("D49").Value = "2.005.60"
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.09"
$ws.Range("E51").Value = "  +0.39%  "
